$d = $word.ActiveDocument

# OOXML package wrapper template used with Range.InsertXML so we can
# control the exact WordprocessingML emitted for each paragraph (proof
# errors, bookmarks, run boundaries) the same way a real edit in Word
# would leave them.
$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Replace-ParagraphXml($paragraph, [string]$innerXml) {
    $r = $paragraph.Range
    # Trim the trailing paragraph mark off the range so we only clear /
    # replace the paragraph's content, keeping the paragraph mark (and
    # therefore the paragraph itself) intact.
    [void]$r.MoveEnd(1, -1)
    $r.Text = ""
    $xml = $pkgOpen + "<w:p>" + $innerXml + "</w:p>" + $pkgClose
    [void]$r.InsertXML($xml)
}

# Find the two target paragraphs by their current text so the script
# isn't dependent on brittle absolute paragraph indices.
$benPara = $null
$againPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "Ben changing things up!") {
        $benPara = $p
    } elseif ($t -eq "Again!") {
        $againPara = $p
    }
}

if ($benPara -ne $null) {
    $innerXml = '<w:proofErr w:type="spellStart"/><w:r><w:t>Labda</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> was changing thing up here!!!!!</w:t></w:r>'
    Replace-ParagraphXml $benPara $innerXml
}

if ($againPara -ne $null) {
    $innerXml = '<w:r><w:t>A</w:t></w:r><w:r><w:t>gain here!!</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
    Replace-ParagraphXml $againPara $innerXml
}
